# SSDM-8226 Updates covered criteria document.
#
# The "Sheet3" tab lists AbstractFieldSearchCriteria subtype rows. A number
# of rows representing criteria that are now "covered" (marked with a "+"
# in column F) get the bold + green highlight applied to the whole row,
# and a handful of rows that previously lacked the "+" marker gain it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Rows that need the "+" marker written into column F (they didn't have
# one before). These are exactly the rows that become fully highlighted
# (bold font on a green fill) in the revised document.
$newPlusRows = @(49,50,51,53,54,63,64,65,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93)

foreach ($r in $newPlusRows) {
    $ws.Cells.Item($r, 6).Value = "+"
}

# Rows that should show the bold-on-green "covered" styling for the full
# row width (columns A:F). This is the set above plus four rows (119, 120,
# 122, 123) that already carried the "+" marker but still need the
# highlight applied. Row 121 intentionally stays unstyled.
$highlightRows = @(49,50,51,53,54,63,64,65,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,119,120,122,123)

foreach ($r in $highlightRows) {
    $rng = $ws.Range("A" + $r + ":F" + $r)
    $rng.Font.Bold = $true
}

# Reflect the new scroll position / active selection recorded for the
# sheet after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 98
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E121").Select()
